# Apply cryptos list price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.991.41"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "1.635.90"
$ws.Range("E3").Value = "  +2.24%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.48%  "

$ws.Range("E6").Value = "  +1.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.61"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.25%  "

$ws.Range("E9").Value = "  +4.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0614"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.37%  "

$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").Value = "1.871.90"
$ws.Range("E12").Value = "  +2.35%  "

$ws.Range("D13").Value = "1.643.31"
$ws.Range("E13").Value = "  +2.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.577"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +25.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.54%  "

$ws.Range("D17").Value = "30.032.65"
$ws.Range("E17").Value = "  +1.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.87%  "

$ws.Range("D20").Value = "0.0₃0709"
$ws.Range("E20").Value = "  +2.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("E22").Value = "  +5.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.81%  "

$ws.Range("E26").Value = "  +2.54%  "

$ws.Range("E27").Value = "  +2.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.26%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  +2.98%  "

$ws.Range("E31").Value = "  +6.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.54%  "

$ws.Range("E33").Value = "  +1.74%  "

$ws.Range("D34").Value = "1.437.74"
$ws.Range("E34").Value = "  +1.10%  "

$ws.Range("E35").Value = "  +8.08%  "

$ws.Range("E36").Value = "  +1.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0172"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.88%  "

$ws.Range("E41").Value = "  +2.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.838"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.63%  "

$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("E44").Value = "  +0.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "55.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.56%  "

$ws.Range("E46").Value = "  +4.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "

$ws.Range("D49").Value = "1.778.66"
$ws.Range("E49").Value = "  +2.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.82%  "

$ws.Range("E51").Value = "  +5.39%  "

